# Add a new localization row (row 6) with key "t1" and its zh_cn / en_us
# translations, matching the Localization.xlsx sheet layout (key, zh_cn, en_us).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

$ws.Range("A6").Value = "t1"
$ws.Range("B6").Value = '<link="glitch">是故障！~!</link>'
$ws.Range("C6").Value = '<link="wave">wave, wave, wave, alright</link>'

# New column (D) width added alongside the existing B/C columns.
$ws.Columns.Item(4).ColumnWidth = 34.1

# Move/refresh selection to the newly added row, like the author's last edit.
[void]$ws.Range("B6").Select()
